# Rename the three "Sector" labels (aligning with Baseline naming) and
# rotate the EU27+UK (column E) values in rows 5-7 so that each data
# series stays attached to its corresponding (renamed) sector label.
#
# Mapping of labels (old -> new), same row position:
#   Row5: "Offshore wind" -> "Onshore wind plants"
#   Row6: "Onshore wind"  -> "Photovoltaic plants"
#   Row7: "PV"            -> "Offshore wind plants"
#
# Because the data must keep following its original sector, the three
# values in column E (rows 5-7) are rotated: new E5 = old E6,
# new E6 = old E7, new E7 = old E5.

$wb = $excel.ActiveWorkbook
$count = $wb.Worksheets.Count

for ($i = 1; $i -le $count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # --- relabel the Sector column (C5:C7) ---
    $ws.Range("C5").Value = "Onshore wind plants"
    $ws.Range("C6").Value = "Photovoltaic plants"
    $ws.Range("C7").Value = "Offshore wind plants"

    # --- rotate the EU27+UK values (column E) so data stays with its sector ---
    $e5 = $ws.Range("E5").Value2
    $e6 = $ws.Range("E6").Value2
    $e7 = $ws.Range("E7").Value2

    $ws.Range("E5").Value = $e6
    $ws.Range("E6").Value = $e7
    $ws.Range("E7").Value = $e5
}
